# Apply updated cryptocurrency price/volume data to Sheet1.
# Values are written through NumberFormat="@" + ClearFormats() so that
# plain-number-looking strings (e.g. "322.75") stay text cells exactly
# like the source data, instead of being auto-coerced to numeric by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" "47.748.01"
Set-TextValue "E2" "  +0.79%  "
Set-TextValue "D3" "2.493.61"
Set-TextValue "E3" "  -0.14%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "322.75"
Set-TextValue "E5" "  -0.11%  "
Set-TextValue "D6" "108.99"
Set-TextValue "E6" "  +0.76%  "
Set-TextValue "E7" "  -0.77%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "E9" "  +1.27%  "
Set-TextValue "D10" "40.63"
Set-TextValue "E10" "  +6.49%  "
Set-TextValue "D11" "0.0814"
Set-TextValue "E11" "  +0.02%  "
Set-TextValue "E12" "  +0.57%  "
Set-TextValue "D13" "18.71"
Set-TextValue "E13" "  +1.37%  "
Set-TextValue "D14" "7.22"
Set-TextValue "E14" "  +0.03%  "
Set-TextValue "D15" "2.884.45"
Set-TextValue "E15" "  -0.08%  "
Set-TextValue "D16" "2.488.79"
Set-TextValue "E16" "  -0.48%  "
Set-TextValue "D17" "0.853"
Set-TextValue "E17" "  -0.20%  "
Set-TextValue "D18" "47.651.60"
Set-TextValue "E18" "  +0.72%  "
Set-TextValue "D19" "13.23"
Set-TextValue "E19" "  +2.31%  "
Set-TextValue "E20" "  -0.82%  "
Set-TextValue "B21" "ShibaInu"
Set-TextValue "C21" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D21" "0.0₃0944"
Set-TextValue "E21" "  +0.18%  "
Set-TextValue "B22" "ImmutableX"
Set-TextValue "C22" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D22" "2.78"
Set-TextValue "E22" "  +13.50%  "
Set-TextValue "D23" "70.80"
Set-TextValue "E23" "  +0.10%  "
Set-TextValue "D24" "247.44"
Set-TextValue "E24" "  -1.45%  "
Set-TextValue "E25" "  -1.26%  "
Set-TextValue "E26" "  +0.09%  "
Set-TextValue "D27" "25.86"
Set-TextValue "E27" "  -1.51%  "
Set-TextValue "D28" "9.99"
Set-TextValue "E28" "  -0.92%  "
Set-TextValue "D29" "0.140"
Set-TextValue "E29" "  +0.97%  "
Set-TextValue "D30" "35.10"
Set-TextValue "E30" "  -0.10%  "
Set-TextValue "E31" "  -0.39%  "
Set-TextValue "D32" "49.78"
Set-TextValue "E32" "  +0.64%  "
Set-TextValue "D33" "19.97"
Set-TextValue "E33" "  +1.02%  "
Set-TextValue "E34" "  -2.69%  "
Set-TextValue "D35" "0.0793"
Set-TextValue "E35" "  -0.34%  "
Set-TextValue "E36" "  +0.14%  "
Set-TextValue "E37" "  -1.54%  "
Set-TextValue "D38" "4.67"
Set-TextValue "E38" "  -0.51%  "
Set-TextValue "E39" "  -1.15%  "
Set-TextValue "D40" "22.59"
Set-TextValue "E40" "  +6.81%  "
Set-TextValue "E41" "  -0.32%  "
Set-TextValue "E42" "  -1.07%  "
Set-TextValue "D43" "119.37"
Set-TextValue "D44" "0.0298"
Set-TextValue "E44" "  -0.23%  "
Set-TextValue "D45" "2.002.09"
Set-TextValue "E45" "  +1.73%  "
Set-TextValue "E46" "  +0.93%  "
Set-TextValue "E47" "  -3.51%  "
Set-TextValue "D48" "1.82"
Set-TextValue "E48" "  +0.89%  "
Set-TextValue "E49" "  -0.57%  "
Set-TextValue "D50" "5.16"
Set-TextValue "E50" "  -2.42%  "
Set-TextValue "D51" "56.92"
Set-TextValue "E51" "  +3.17%  "
